$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 10.65
$ws.Range("C4").Value = 9.35
$ws.Range("E4").Value = 9.9
$ws.Range("D5").Value = 10.1
$ws.Range("F5").Value = 10.28
$ws.Range("G5").Value = 9.32
$ws.Range("H5").Value = 8.35
$ws.Range("E6").Value = 9.720000000000001
$ws.Range("G6").Value = 10.18
$ws.Range("E7").Value = 10.68
$ws.Range("F7").Value = 9.82
$ws.Range("E8").Value = 11.65
